$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Statut" column (H) values to reflect current order status
$ws.Range("H2").Value = "Expédié"
$ws.Range("H3").Value = "Expédié"
$ws.Range("H4").Value = "Expédié"
$ws.Range("H5").Value = "Stock insuffisant"
$ws.Range("H6").Value = "Expédié"

# Reset the view: scroll back to A1 and set zoom to 85%
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$win.Zoom = 85

$ws.Select()
$ws.Range("H7").Select()
